# Fixed small bugs according to Trello
#
# 1. The "_GoBack" bookmark that used to sit at the end of the
#    "Do the same for the inner walls and bottom wall" paragraph is
#    removed from there.
# 2. Three new paragraphs are appended at the end of the document:
#      - an empty paragraph
#      - a paragraph "Input – 5 x 4" (now carrying the _GoBack bookmark
#        start)
#      - a paragraph with the pseudo-code for-loop (carrying the
#        _GoBack bookmark end)

$d = $word.ActiveDocument

# --- Step 1: remove the old _GoBack bookmark -------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# --- Step 2: append the new paragraphs at the end of the document ----
$newParagraphsXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:r><w:t>Input &#8211; 5 x 4</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="gramStart"/><w:r><w:t>for(</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramEnd"/><w:r><w:t>int</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">=1; </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>i</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">&lt;5; </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>i</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>++)</w:t></w:r><w:r><w:br/></w:r><w:r><w:tab/><w:t xml:space="preserve">//Make the four vertices where </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>i</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>/total length gives the points</w:t></w:r><w:r><w:br/></w:r><w:r><w:tab/><w:t>//Connect them, fill face, separate</w:t></w:r><w:r><w:br/><w:t>//Select all, join</w:t></w:r><w:bookmarkEnd w:id="0"/></w:p>
'@

$endRange = $d.Range($d.Content.End, $d.Content.End)
[void]$endRange.InsertXML($newParagraphsXml)
